# Adds header labels to column A on each "chart data" sheet, fixes
# accentuation on several labels, removes the bold/bordered header style
# from the now-plain label cells, drops the obsolete "Teto" row from the
# emissions sheet, and refreshes the cost-sheet header/values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the exact same row layout (Fonte/Tecnologia table).
# ---------------------------------------------------------------------
$techSheets = @(
    $wb.Worksheets.Item(1),
    $wb.Worksheets.Item(2),
    $wb.Worksheets.Item(3),
    $wb.Worksheets.Item(4)
)

foreach ($ws in $techSheets) {
    # New header for column A, styled like the year headers in B1:E1.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accentuation on a few labels.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # Remove the bold/bordered formatting that used to sit on A2:A12
    # (only the new A1 header keeps it now).
    $ws.Range("A2:A12").ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet 5 - "Emissoes Totais (MtCO2eq)": new header, fixed labels,
# styling removed from data rows, and the "Teto" row is dropped.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").ClearFormats()

$ws5.Rows("4:4").Delete()

# ---------------------------------------------------------------------
# Sheet 6 - "Custo Total (bilhões de R$)": new header, renamed/fixed
# labels, styling removed from data rows, and refreshed values.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws1 = $wb.Worksheets.Item(1)

# B1 becomes a "2015" year label: copy both value + format from sheet 1's
# B1 so the text stays a real string (not auto-coerced to a number).
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4104)

# A1 gets a new "Tipo Expansão" header: copy B1's format only, then set
# the (different) text explicitly.
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 587
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
$ws6.Range("A2:A3").ClearFormats()
